$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wording of existing test-case descriptions ("whether" -> "that")
$ws.Range("C12").Value = "This is to test that users are not able to add a To Do Item with no input"
$ws.Range("C13").Value = "This is to test that users can add a To Do Item with valid parameters"

# Add the new "test_AddExistingItem" test case as row 14
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "test_AddExistingItem"
$ws.Range("C14").Value = "This is to test that users can add a To Do item that already exists"
$ws.Range("E14").Value = "Item is added to the list"

# Move the active selection to reflect where the author ended up (E15)
$ws.Range("E15").Select()
